$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.380.65"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "2.647.41"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'596.24"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").Value = "'158.85"
$ws.Range("E6").Value = "  +2.74%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("D9").Value = "2.646.83"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("E10").Value = "  -2.31%  "
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D12").Value = "'5.27"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").Value = "'0.352"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").Value = "3.131.56"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "'0.0000187"
$ws.Range("E16").Value = "  -3.22%  "
$ws.Range("D17").Value = "68.291.41"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "2.616.41"
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("D19").Value = "'11.57"
$ws.Range("E19").Value = "  +1.88%  "
$ws.Range("D20").Value = "'363.64"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "'4.40"
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("E23").Value = "  -1.82%  "
$ws.Range("D24").Value = "'2.09"
$ws.Range("E24").Value = "  +1.16%  "
$ws.Range("D25").Value = "'74.83"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").Value = "'9.91"
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("D28").Value = "2.796.93"
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("E29").Value = "  -2.75%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").Value = "'565.67"
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("D32").Value = "'8.04"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").Value = "'1.64"
$ws.Range("E35").Value = "  +4.44%  "
$ws.Range("E36").Value = "  -1.36%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "'160.56"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").Value = "'19.67"
$ws.Range("E39").Value = "  +1.89%  "
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("E41").Value = "  -0.94%  "
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("D43").Value = "'2.64"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").Value = "0.0₆0320"
$ws.Range("E44").Value = "  -5.59%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "'158.31"
$ws.Range("E46").Value = "  +1.42%  "
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("D48").Value = "'21.86"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").Value = "'0.0779"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("E51").Value = "  +1.99%  "
